$wb = $excel.ActiveWorkbook

# OFF sheet - Row 3 ("R" / Road) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 321
$wsOff.Range("C3").Value = 206
$wsOff.Range("D3").Value = 73
$wsOff.Range("E3").Value = 37

# DEF sheet - Row 3 ("R" / Road) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 399
$wsDef.Range("C3").Value = 276
$wsDef.Range("D3").Value = 86
$wsDef.Range("E3").Value = 39
